$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I3").Value = 1.400202336862439
$ws.Range("J3").Value = 0.4842730983059223
$ws.Range("K3").Value = -0.5938596744236889
$ws.Range("L3").Value = 2.059269257706759
